# Update the "length" column (R) to hold a data-type string, and clear the
# "sampling" column (H) for DIGITAL_OUTPUT / ANALOG_OUTPUT rows.
#
# Rules (derived from the target diff):
#   - rows whose "length" (R) starts with "BIT "                -> R = "1bit"
#   - rows whose system_category (E) is "COMMAND"                -> R = "1bit"
#   - rows whose system_category (E) is "ANALOG_OUTPUT"          -> R = "s16"
#   - all remaining data rows (CONFIG_PARAMETER, SET_POINT,
#     DIGITAL_OUTPUT, STATUS, ...)                                -> R = "f32cdab"
#   - rows whose system_category (E) is "DIGITAL_OUTPUT" or
#     "ANALOG_OUTPUT" additionally get their sampling (H) reset
#     from 60 to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $category = $ws.Cells.Item($r, 5).Value()
    $length = $ws.Cells.Item($r, 18).Value()
    $lengthText = [string]$length

    if ($lengthText.StartsWith("BIT")) {
        $newLength = "1bit"
    }
    elseif ($category -eq "COMMAND") {
        $newLength = "1bit"
    }
    elseif ($category -eq "ANALOG_OUTPUT") {
        $newLength = "s16"
    }
    else {
        $newLength = "f32cdab"
    }

    $ws.Cells.Item($r, 18).Value = $newLength

    if ($category -eq "DIGITAL_OUTPUT" -or $category -eq "ANALOG_OUTPUT") {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
